$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.85

# Row 6
$ws.Range("G6").Value = 1.31
$ws.Range("H6").Value = 4.9
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 1.75
$ws.Range("K6").Value = 2.6
$ws.Range("N6").Value = 9.5
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4.8
$ws.Range("Q6").Value = 1.45
$ws.Range("R6").Value = 2.55
$ws.Range("S6").Value = 1.26
$ws.Range("T6").Value = 3.5
$ws.Range("U6").Value = 1.75
$ws.Range("V6").Value = 1.95
$ws.Range("W6").Value = 9.25
$ws.Range("X6").Value = 7.5
$ws.Range("AA6").Value = 10.25
$ws.Range("AB6").Value = 22
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 10.5
$ws.Range("AE6").Value = 18.5
$ws.Range("AF6").Value = 70
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 60
$ws.Range("AJ6").Value = 200
$ws.Range("AL6").Value = 60
$ws.Range("AM6").Value = 450
$ws.Range("AN6").Value = 3.35
$ws.Range("AO6").Value = 5.8
$ws.Range("AP6").Value = 14
$ws.Range("AQ6").Value = 14.5
$ws.Range("AR6").Value = 35
$ws.Range("AS6").Value = 150
$ws.Range("AT6").Value = 3.5
$ws.Range("AV6").Value = 60
$ws.Range("AW6").Value = 9
$ws.Range("AY6").Value = 37
$ws.Range("AZ6").Value = 250
$ws.Range("BB6").Value = 400

# Row 7
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 3.65
$ws.Range("I7").Value = 1.6
$ws.Range("J7").Value = 5.3
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 2.2
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 7.6
$ws.Range("O7").Value = 1.27
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 1.82
$ws.Range("R7").Value = 1.91
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.88
$ws.Range("W7").Value = 14
$ws.Range("X7").Value = 30
$ws.Range("Y7").Value = 16.5
$ws.Range("Z7").Value = 100
$ws.Range("AA7").Value = 55
$ws.Range("AB7").Value = 55
$ws.Range("AC7").Value = 7.6
$ws.Range("AD7").Value = 7.3
$ws.Range("AE7").Value = 16
$ws.Range("AF7").Value = 75
$ws.Range("AG7").Value = 6.9
$ws.Range("AH7").Value = 7.6
$ws.Range("AJ7").Value = 12
$ws.Range("AK7").Value = 12.5
$ws.Range("AL7").Value = 25
$ws.Range("AM7").Value = 600
$ws.Range("AN7").Value = 6.7
$ws.Range("AO7").Value = 30
$ws.Range("AU7").Value = 7.8
$ws.Range("AW7").Value = 3.4
$ws.Range("AX7").Value = 8
$ws.Range("AZ7").Value = 27
$ws.Range("BA7").Value = 60
$ws.Range("BB7").Value = 250
